$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $text) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t -eq $text -or $t.TrimEnd([char]13) -eq $text) {
            return $p
        }
    }
    return $null
}

function Get-ParagraphIndexByStart($doc, $startPos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Start -eq $startPos) {
            return $i
        }
    }
    return -1
}

# Drop a fully-bold (w:b + w:bCs, on the run AND the paragraph mark) copy of
# $formattedText right after $afterPara, then retarget its visible text to
# $newText. Returns the resulting paragraph.
#
# Directly toggling Range.Font.Bold / Range.Font.BoldBi on a paragraph only
# ever stamps the "ascii" (w:b) flag onto the paragraph-mark run properties,
# never the complex-script (w:bCs) one (it does manage to set both on the
# text run itself), so instead we clone a paragraph's FormattedText that
# already carries w:b + w:bCs in both places and reuse it as a donor.
#
# Inserting FormattedText exactly at the end of the document body merges it
# into the preceding paragraph instead of splitting off a new one, so when
# $afterPara is the last paragraph in the document we temporarily add a
# throwaway paragraph after it first to guarantee there's a following
# paragraph boundary for the split to land on.
function Insert-BoldParagraphAfter($doc, $afterPara, $formattedText, $newText) {
    $isLastParagraph = ($afterPara.Range.End -eq $doc.Content.End)
    if ($isLastParagraph) {
        $endPos = $afterPara.Range.End
        $doc.Range($endPos, $endPos).InsertParagraphAfter() | Out-Null
    }

    $insertPos = $afterPara.Range.End
    $doc.Range($insertPos, $insertPos).FormattedText = $formattedText

    $afterIndex = Get-ParagraphIndexByStart $doc $afterPara.Range.Start
    $newPara = $doc.Paragraphs.Item($afterIndex + 1)
    $doc.Range($newPara.Range.Start, $newPara.Range.End - 1).Text = $newText
    return $newPara
}

# ---------------------------------------------------------------------------
# 1. Make the "Clone Repository" heading bold, matching the style of the
#    other section headings in the document (w:b + w:bCs on both the
#    paragraph mark and the run).
# ---------------------------------------------------------------------------
$donorPara = Get-ParagraphByText $d "Tell Git who you are"
$donorFormattedText = $donorPara.Range.FormattedText

$oldCloneHeading = Get-ParagraphByText $d "Clone Repository"
$cmdAfterHeading = Get-ParagraphIndexByStart $d $oldCloneHeading.Range.Start
$precedingPara = $d.Paragraphs.Item($cmdAfterHeading - 1)

Insert-BoldParagraphAfter $d $oldCloneHeading $donorFormattedText "Clone Repository" | Out-Null

# Remove the old, plain "Clone Repository" paragraph.
$oldCloneHeading.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Append the new "Pull" section: a bold "Pull " heading followed by a
#    "git pull origin master" list item that shares the numbering list
#    (numId 8) used by the "git clone <https link>" item above it.
# ---------------------------------------------------------------------------
$cloneCmdPara = Get-ParagraphByText $d "git clone <https link>"

$pullHeadingPara = Insert-BoldParagraphAfter $d $cloneCmdPara $donorFormattedText "Pull "

# The helper above leaves a placeholder paragraph behind whenever it had to
# extend the document first (true here, since "git clone <https link>" is
# the last paragraph) — that placeholder already inherited the preceding
# ListParagraph / numId 8 formatting, so just fill it in as the new list
# item instead of inserting yet another paragraph.
$pullHeadingIndex = Get-ParagraphIndexByStart $d $pullHeadingPara.Range.Start
$gitPullPara = $d.Paragraphs.Item($pullHeadingIndex + 1)
$isPlaceholder = ($gitPullPara.Range.End - $gitPullPara.Range.Start) -eq 1
if (-not $isPlaceholder) {
    $endPos = $pullHeadingPara.Range.End
    $d.Range($endPos, $endPos).InsertParagraphAfter() | Out-Null
    $gitPullPara = $d.Paragraphs.Item($pullHeadingIndex + 1)
}
$d.Range($gitPullPara.Range.Start, $gitPullPara.Range.End - 1).Text = "git pull origin master"

Write-Output "Edit complete"
